$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), reusing the same header style as H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Fill I2:I33 with constant 1, and J2:J33 with the same values as H2:H33
$lastRow = 33
for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
